# Auto-generated script applying cryptos.xlsx price/volume updates
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.306.99'
$ws.Range('E2').Value = '  +1.07%  '
$ws.Range('D3').Value = '1.619.45'
$ws.Range('E3').Value = '  +1.75%  '
$ws.Range('D5').Value = '''212.00'
$ws.Range('D5').Style = "Normal"
$ws.Range('E5').Value = '  +0.55%  '
$ws.Range('E6').Value = '  -0.02%  '
$ws.Range('D7').Value = '''0.483'
$ws.Range('D7').Style = "Normal"
$ws.Range('E7').Value = '  +0.53%  '
$ws.Range('E8').Value = '  +0.27%  '
$ws.Range('D10').Value = '''18.74'
$ws.Range('D10').Style = "Normal"
$ws.Range('E10').Value = '  +4.24%  '
$ws.Range('D11').Value = '''0.0814'
$ws.Range('D11').Style = "Normal"
$ws.Range('E11').Value = '  +0.77%  '
$ws.Range('D12').Value = '1.845.43'
$ws.Range('E12').Value = '  +1.75%  '
$ws.Range('D13').Value = '1.624.93'
$ws.Range('E13').Value = '  +2.04%  '
$ws.Range('D14').Value = '''4.03'
$ws.Range('D14').Style = "Normal"
$ws.Range('E14').Value = '  +0.97%  '
$ws.Range('D15').Value = '''0.517'
$ws.Range('D15').Style = "Normal"
$ws.Range('E15').Value = '  +1.19%  '
$ws.Range('D16').Value = '26.308.46'
$ws.Range('E16').Value = '  +1.10%  '
$ws.Range('D17').Value = '''62.31'
$ws.Range('D17').Style = "Normal"
$ws.Range('E17').Value = '  +3.47%  '
$ws.Range('D18').Value = '0.0₃0726'
$ws.Range('E18').Value = '  +0.43%  '
$ws.Range('E19').Value = '  -0.04%  '
$ws.Range('D20').Value = '''201.60'
$ws.Range('D20').Style = "Normal"
$ws.Range('E20').Value = '  -0.05%  '
$ws.Range('D21').Value = '''4.27'
$ws.Range('D21').Style = "Normal"
$ws.Range('E21').Value = '  +0.88%  '
$ws.Range('D22').Value = '''9.32'
$ws.Range('D22').Style = "Normal"
$ws.Range('E22').Value = '  +1.34%  '
$ws.Range('D23').Value = '''6.03'
$ws.Range('D23').Style = "Normal"
$ws.Range('E23').Value = '  +0.66%  '
$ws.Range('E24').Value = '  -5.12%  '
$ws.Range('D25').Value = '''144.53'
$ws.Range('D25').Style = "Normal"
$ws.Range('E25').Value = '  +0.65%  '
$ws.Range('E26').Value = '  -0.10%  '
$ws.Range('E27').Value = '  -1.23%  '
$ws.Range('E28').Value = '  +0.26%  '
$ws.Range('E29').Value = '  +1.29%  '
$ws.Range('D30').Value = '''0.0516'
$ws.Range('D30').Style = "Normal"
$ws.Range('E30').Value = '  +8.71%  '
$ws.Range('E31').Value = '  +0.74%  '
$ws.Range('E32').Value = '  +1.84%  '
$ws.Range('E33').Value = '  -0.47%  '
$ws.Range('E34').Value = '  +1.26%  '
$ws.Range('E35').Value = '  +2.53%  '
$ws.Range('D36').Value = '1.176.89'
$ws.Range('E36').Value = '  +4.25%  '
$ws.Range('E37').Value = '  +0.62%  '
$ws.Range('D38').Value = '''0.808'
$ws.Range('D38').Style = "Normal"
$ws.Range('E38').Value = '  +2.10%  '
$ws.Range('E39').Value = '  -0.01%  '
$ws.Range('E40').Value = '  +0.12%  '
$ws.Range('E41').Value = '  +1.19%  '
$ws.Range('E42').Value = '  +4.20%  '
$ws.Range('D43').Value = '''0.783'
$ws.Range('D43').Style = "Normal"
$ws.Range('E43').Value = '  +0.52%  '
$ws.Range('D44').Value = '1.757.63'
$ws.Range('E44').Value = '  +1.94%  '
$ws.Range('D45').Value = '''92.52'
$ws.Range('D45').Style = "Normal"
$ws.Range('E45').Value = '  +0.36%  '
$ws.Range('E46').Value = '  +2.81%  '
$ws.Range('D47').Value = '''53.74'
$ws.Range('D47').Style = "Normal"
$ws.Range('E47').Value = '  +0.10%  '
$ws.Range('E48').Value = '  +1.05%  '
$ws.Range('D49').Value = '''0.409'
$ws.Range('D49').Style = "Normal"
$ws.Range('E49').Value = '  +0.51%  '
$ws.Range('E50').Value = '  -0.21%  '
$ws.Range('E51').Value = '  +2.05%  '
